$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell $ws "D2" "29.001.98"
Set-TextCell $ws "E2" "  +0.38%  "
Set-TextCell $ws "D3" "1.887.04"
Set-TextCell $ws "E3" "  -1.13%  "
Set-TextCell $ws "E4" "  +0.15%  "
Set-TextCell $ws "D5" "330.53"
Set-TextCell $ws "E5" "  -2.31%  "
Set-TextCell $ws "D6" "1.0000"
Set-TextCell $ws "D7" "0.4590"
Set-TextCell $ws "E7" "  -2.78%  "
Set-TextCell $ws "D8" "0.4058"
Set-TextCell $ws "E8" "  +0.39%  "
Set-TextCell $ws "E9" "  -0.50%  "
Set-TextCell $ws "E10" "  -2.41%  "
Set-TextCell $ws "D11" "0.9917"
Set-TextCell $ws "E11" "  -2.85%  "
Set-TextCell $ws "E12" "  -3.53%  "
Set-TextCell $ws "D13" "1.879.06"
Set-TextCell $ws "E13" "  -0.99%  "
Set-TextCell $ws "D14" "5.906"
Set-TextCell $ws "E14" "  -2.99%  "
Set-TextCell $ws "D15" "7.067"
Set-TextCell $ws "E15" "  -3.85%  "
Set-TextCell $ws "E16" "  +0.10%  "
Set-TextCell $ws "D17" "88.35"
Set-TextCell $ws "E17" "  -3.21%  "
Set-TextCell $ws "D18" "0.00001028"
Set-TextCell $ws "E18" "  -2.30%  "
Set-TextCell $ws "D19" "0.06545"
Set-TextCell $ws "E19" "  -1.00%  "
Set-TextCell $ws "D20" "17.43"
Set-TextCell $ws "E20" "  -2.03%  "
Set-TextCell $ws "E21" "  +0.04%  "
Set-TextCell $ws "D22" "29.006.66"
Set-TextCell $ws "E22" "  +0.30%  "
Set-TextCell $ws "D23" "5.416"
Set-TextCell $ws "E23" "  -2.49%  "
Set-TextCell $ws "D24" "11.33"
Set-TextCell $ws "E24" "  +1.53%  "
Set-TextCell $ws "D25" "2.207"
Set-TextCell $ws "E25" "  -2.68%  "
Set-TextCell $ws "D26" "2.103.57"
Set-TextCell $ws "E26" "  -0.61%  "
Set-TextCell $ws "D27" "156.53"
Set-TextCell $ws "E27" "  -2.88%  "
Set-TextCell $ws "D28" "19.58"
Set-TextCell $ws "E28" "  -2.04%  "
Set-TextCell $ws "D29" "2.105"
Set-TextCell $ws "E29" "  -2.93%  "
Set-TextCell $ws "D30" "5.417"
Set-TextCell $ws "E30" "  -2.08%  "
Set-TextCell $ws "D31" "117.76"
Set-TextCell $ws "E31" "  -2.32%  "
Set-TextCell $ws "D32" "1.003"
Set-TextCell $ws "E32" "  -0.86%  "
Set-TextCell $ws "D33" "0.09331"
Set-TextCell $ws "E33" "  -2.53%  "
Set-TextCell $ws "E34" "  -1.54%  "
Set-TextCell $ws "D35" "1.410"
Set-TextCell $ws "E35" "  -0.34%  "
Set-TextCell $ws "D36" "5.280"
Set-TextCell $ws "E36" "  -2.31%  "
Set-TextCell $ws "D37" "0.06053"
Set-TextCell $ws "E37" "  -2.40%  "
Set-TextCell $ws "D38" "0.02217"
Set-TextCell $ws "E38" "  -3.04%  "
Set-TextCell $ws "D39" "8.265"
Set-TextCell $ws "E39" "  -4.76%  "
Set-TextCell $ws "D40" "1.179"
Set-TextCell $ws "E40" "  -1.53%  "
Set-TextCell $ws "D41" "0.9999"
Set-TextCell $ws "E41" "  +0.02%  "
Set-TextCell $ws "D42" "0.5783"
Set-TextCell $ws "E42" "  -4.00%  "
Set-TextCell $ws "E43" "  -4.03%  "
Set-TextCell $ws "E44" "  -4.29%  "
Set-TextCell $ws "D45" "1.259"
Set-TextCell $ws "E45" "  -1.93%  "
Set-TextCell $ws "D46" "0.07530"
Set-TextCell $ws "E46" "  +1.93%  "
Set-TextCell $ws "D47" "2.278"
Set-TextCell $ws "E47" "  +7.09%  "
Set-TextCell $ws "D48" "11.98"
Set-TextCell $ws "E48" "  -3.17%  "
Set-TextCell $ws "D49" "0.5455"
Set-TextCell $ws "E49" "  -3.24%  "
Set-TextCell $ws "D50" "1.897"
Set-TextCell $ws "E50" "  -4.24%  "

# Row 51: coin renamed from Elrond to Quant
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws "D51" "111.09"
Set-TextCell $ws "E51" "  -1.83%  "
